$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.952.91"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.557.01"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.07"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.02%  "

$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0589"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.05%  "

$ws.Range("D12").Value = "1.778.91"
$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("D13").Value = "1.556.57"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("E15").Value = "  +1.70%  "

$ws.Range("D16").Value = "26.951.47"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.80"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("D19").Value = "0.0₃0696"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("E20").Value = "  +1.25%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.18%  "

$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.46%  "

$ws.Range("E28").Value = "  +0.50%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0470"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("E31").Value = "  -1.18%  "

$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("D33").Value = "1.422.67"
$ws.Range("E33").Value = "  +4.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.08"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.13%  "

$ws.Range("E35").Value = "  +3.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.978"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.81%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("E43").Value = "  +3.46%  "

$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.92%  "

$ws.Range("E46").Value = "  +1.68%  "

$ws.Range("D47").Value = "1.692.80"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.56"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("D50").Value = "0.0₇0990"
$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.30%  "
